$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrangeHrmData")
$ws.Activate()

# Change B2:B6 from text last names to numeric values
$ws.Range("B2").Value = 123
$ws.Range("B3").Value = 456
$ws.Range("B4").Value = 789
$ws.Range("B5").Value = 123
$ws.Range("B6").Value = 345

# Apply thick bottom border + row height to rows 2-6 (matching row 1/6 previously)
$rng = $ws.Range("A2:D6")
$rng.Borders.Item(9).LineStyle = 1
$rng.Borders.Item(9).Weight = -4138

$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75

# View changes: scroll & selection
$ws.Range("B5").Select()

$ws.PageSetup.Orientation = 1
